$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -1        # Total P&L %
$summary.Range("B6").Value = 50        # Total Trades
$summary.Range("B9").Value = 38        # Win Rate %

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 50         # MarketMaking Trades
$status.Range("G4").Value = 38         # MarketMaking Win Rate %

# --- New trade row (#50) appended to "All Trades" and "MarketMaking" sheets ---
$sheetsToUpdate = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetsToUpdate) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Force Date/Time columns to be stored as plain text (matches existing rows),
    # avoiding Excel's automatic date/time serial-number conversion.
    $ws.Range("B51:C51").NumberFormat = "@"

    $ws.Range("A51").Value = 50
    $ws.Range("B51").Value = "2026-02-17"
    $ws.Range("C51").Value = "13:28:35"
    $ws.Range("D51").Value = "MarketMaking"
    $ws.Range("E51").Value = "UP"
    $ws.Range("F51").Value = 0.98
    $ws.Range("G51").Value = 0.98
    $ws.Range("H51").Value = "CLOSED"
    $ws.Range("I51").Value = 0
    $ws.Range("J51").Value = 0
    $ws.Range("K51").Value = 97.5
    $ws.Range("L51").Value = 0
    $ws.Range("M51").Value = 0
    $ws.Range("N51").Value = 0.6
    $ws.Range("O51").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P51").Value = "early_exit"
    $ws.Range("Q51").Value = 0.14

    # Restore the default "Normal" style so the text cells don't retain the
    # custom number-format style index created above.
    $ws.Range("B51:C51").Style = "Normal"
}
